# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1179.8334
$ws.Range("I98").Value = 1039.4445
$ws.Range("J98").Value = 1601
$ws.Range("K98").Value = 1039.4445
$ws.Range("L98").Value = 1601
$ws.Range("M98").Value = 458.5554999999999
$ws.Range("N98").Value = -4597

$ws.Range("H112").Value = 13938387
$ws.Range("J112").Value = 15038759
$ws.Range("L112").Value = 45116277
$ws.Range("N112").Value = -45118493

$ws.Range("H122").Value = 1179.8334
$ws.Range("I122").Value = 1039.4445
$ws.Range("J122").Value = 1601
$ws.Range("K122").Value = 3118.3335
$ws.Range("L122").Value = 4803
$ws.Range("M122").Value = -668.3335000000002
$ws.Range("N122").Value = -9703

$ws.Range("H138").Value = 2157.9473
$ws.Range("I138").Value = 1180.6285
$ws.Range("J138").Value = 2728.05
$ws.Range("K138").Value = 3541.8855
$ws.Range("L138").Value = 8184.150000000001
$ws.Range("M138").Value = 1598.1145
$ws.Range("N138").Value = -18464.15

$ws.Range("H141").Value = 2271.8386
$ws.Range("I141").Value = 1292.091
$ws.Range("K141").Value = 3876.273
$ws.Range("M141").Value = 1303.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 6407.6924
$ws.Range("I6").Value = 2000
$ws.Range("J6").Value = 7209.091
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 7209.091
$ws.Range("M6").Value = -1827
$ws.Range("N6").Value = -7555.091

$ws.Range("H32").Value = 3078.85
$ws.Range("I32").Value = 3078.85
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3078.85
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2791.85
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 13507.875
$ws.Range("I45").Value = 15080.429
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 15080.429
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -14703.429
$ws.Range("N45").Value = -3254

$ws.Range("H76").Value = 36500
$ws.Range("J76").Value = 36500
$ws.Range("L76").Value = 36500
$ws.Range("N76").Value = -37176

$ws.Range("H79").Value = 36500
$ws.Range("J79").Value = 36500
$ws.Range("L79").Value = 36500
$ws.Range("N79").Value = -38840

$ws.Range("H92").Value = 25333.166
$ws.Range("J92").Value = 25333.166
$ws.Range("L92").Value = 25333.166
$ws.Range("N92").Value = -30325.166

$ws.Range("H122").Value = 571712.6
$ws.Range("I122").Value = 803104.6
$ws.Range("J122").Value = 2132.4614
$ws.Range("K122").Value = 2409313.8
$ws.Range("L122").Value = 6397.3842
$ws.Range("M122").Value = -2406863.8
$ws.Range("N122").Value = -11297.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 3731.5
$ws.Range("I37").Value = 663
$ws.Range("K37").Value = 663
$ws.Range("M37").Value = -526

$ws.Range("H99").Value = 52633388
$ws.Range("I99").Value = 83334420
$ws.Range("K99").Value = 83334420
$ws.Range("M99").Value = -83332922

$ws.Range("H107").Value = 1410.625
$ws.Range("I107").Value = 1411.3334
$ws.Range("K107").Value = 1411.3334
$ws.Range("M107").Value = 508.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 34500
$ws.Range("J28").Value = 34500
$ws.Range("L28").Value = 34500
$ws.Range("N28").Value = -34990

$ws.Range("H58").Value = 229206.02
$ws.Range("I58").Value = 1627.75
$ws.Range("J58").Value = 359250.75
$ws.Range("K58").Value = 1627.75
$ws.Range("L58").Value = 359250.75
$ws.Range("M58").Value = -1424.75
$ws.Range("N58").Value = -359656.75

$ws.Range("H94").Value = 8060.5557
$ws.Range("I94").Value = 7487
$ws.Range("K94").Value = 7487
$ws.Range("M94").Value = -7036

$ws.Range("H99").Value = 24857.143
$ws.Range("I99").Value = 28500
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 28500
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -27002
$ws.Range("N99").Value = -5996

$ws.Range("H107").Value = 991.69696
$ws.Range("I107").Value = 1038.65
$ws.Range("J107").Value = 919.46155
$ws.Range("K107").Value = 1038.65
$ws.Range("L107").Value = 919.46155
$ws.Range("M107").Value = 881.3499999999999
$ws.Range("N107").Value = -4759.46155

$ws.Range("H126").Value = 24857.143
$ws.Range("I126").Value = 28500
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 85500
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -83030
$ws.Range("N126").Value = -13940

$ws.Range("H136").Value = 229206.02
$ws.Range("I136").Value = 1627.75
$ws.Range("J136").Value = 359250.75
$ws.Range("K136").Value = 4883.25
$ws.Range("L136").Value = 1077752.25
$ws.Range("M136").Value = -2333.25
$ws.Range("N136").Value = -1082852.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1130.3334
$ws.Range("I16").Value = 195.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 586.5
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -413.5
$ws.Range("N16").Value = -9346

$ws.Range("H20").Value = 7666.6665
$ws.Range("J20").Value = 7666.6665
$ws.Range("L20").Value = 22999.9995
$ws.Range("N20").Value = -23453.9995

$ws.Range("H26").Value = 747.8182
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 812.6
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 2437.8
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -3013.8

$ws.Range("H29").Value = 749.0833
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 1023.625
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 3070.875
$ws.Range("M29").Value = -323
$ws.Range("N29").Value = -3624.875

$ws.Range("H39").Value = 2107
$ws.Range("J39").Value = 2107
$ws.Range("L39").Value = 6321
$ws.Range("N39").Value = -6909

$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -10068

$ws.Range("H43").Value = 7000
$ws.Range("J43").Value = 7000
$ws.Range("L43").Value = 21000
$ws.Range("N43").Value = -21228

$ws.Range("H131").Value = 1887816.9
$ws.Range("J131").Value = 1081.6818
$ws.Range("L131").Value = 3245.0454
$ws.Range("N131").Value = -13325.0454

$ws.Range("H136").Value = 5383.75
$ws.Range("I136").Value = 7963.5713
$ws.Range("J136").Value = 4523.8096
$ws.Range("K136").Value = 23890.7139
$ws.Range("L136").Value = 13571.4288
$ws.Range("M136").Value = -18790.7139
$ws.Range("N136").Value = -23771.4288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 17000
$ws.Range("J103").Value = 17000
$ws.Range("L103").Value = 17000
$ws.Range("N103").Value = -19344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 24880.334
$ws.Range("J101").Value = 24880.334
$ws.Range("L101").Value = 24880.334
$ws.Range("N101").Value = -31370.334

$ws.Range("H105").Value = 48307.5
$ws.Range("J105").Value = 48307.5
$ws.Range("L105").Value = 48307.5
$ws.Range("N105").Value = -55295.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 7220
$ws.Range("J28").Value = 7220
$ws.Range("L28").Value = 7220
$ws.Range("N28").Value = -7916

$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

$ws.Range("H126").Value = 804.3913
$ws.Range("I126").Value = 711.05
$ws.Range("J126").Value = 1426.6666
$ws.Range("K126").Value = 2133.15
$ws.Range("L126").Value = 4279.9998
$ws.Range("M126").Value = 336.8500000000004
$ws.Range("N126").Value = -9219.9998

$ws.Range("H136").Value = 2451.5273
$ws.Range("I136").Value = 2585.138
$ws.Range("J136").Value = 2302.5
$ws.Range("K136").Value = 7755.414
$ws.Range("L136").Value = 6907.5
$ws.Range("M136").Value = -5205.414
$ws.Range("N136").Value = -12007.5

$ws.Range("H137").Value = 39653.75
$ws.Range("J137").Value = 39653.75
$ws.Range("L137").Value = 39653.75
$ws.Range("N137").Value = -49853.75

$ws.Range("H138").Value = 45229
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 45229
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 45229
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -55509
